function Set-RowValues {
    param($ws, $row, $startCol, $values)
    $col = $startCol
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

function Set-IndexCellStyle {
    # Reproduce the bold/centered/bordered formatting ("s=1" in the original
    # file) used on every A-column HKL-index cell, since newly-used rows
    # start out with no explicit formatting. Copy the format from an
    # existing, already-styled index cell (A9) rather than setting font /
    # alignment / border properties individually.
    param($ws, $row)
    $ws.Cells.Item(9, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook previously held 14 sampling schemes (rows 3-16). New runs were
# computed for three additional "Spiral" sampling schemes, plus a recomputed
# "Gaussian-Quadrature" row. These are placed right after the
# "Ring Perpendicular to TD" row (row 9) - pushing the remaining schemes
# (NoRotation-tilt60deg ... HexGrid-60degTilt5degRes) down by three rows and
# extending the used range through row 19.

# --- Row 10: Gaussian-Quadrature (recomputed) ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
Set-RowValues $ws 10 3 @(1.031533375272873,0.9183273690739626,1.008475206288345,1.031533375272873,0.9473117055122916,1.038744308907424,1.013762454306349,0.9183273690739626,0.9634012876811539,0.9974673314770136,0.9930257365602078)

# --- Row 11: Spiral-90deg-10rot-5space (new) ---
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
Set-RowValues $ws 11 3 @(0.9486757638535821,0.9308341099036397,1.04435713571757,0.9486757638535821,0.9236387522738029,1.165853963576787,1.016468887277038,0.9308341099036397,0.9875956228106048,0.9681356933320935,1.004971435433737)

# --- Row 12: Spiral-90deg-15rot-5space (new) ---
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
Set-RowValues $ws 12 3 @(0.9485100606339691,0.9327026992368069,1.043709298156651,0.9485100606339691,0.924854684485083,1.164019313365911,1.015902408354741,0.9327026992368069,0.9882059986967289,0.968358029665349,1.00494974403886)

# --- Row 13: Spiral-90deg-10rot-3space (new) ---
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
Set-RowValues $ws 13 3 @(0.9484805177774732,0.9314055662718768,1.044229379429764,0.9484805177774732,0.9239232877818482,1.16561918463028,1.016301734374657,0.9314055662718768,0.9878174728508202,0.9681489953141468,1.00499327837765)

# --- Row 14: NoRotation-tilt60deg (existing data, shifted down 4 rows) ---
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
Set-RowValues $ws 14 3 @(0.9542800000000012,0.8961880000000001,1.054419999999999,0.9542800000000012,0.9202959999999999,1.155904000000001,1.027356,0.8961880000000001,0.9753039999999993,0.9647920000000002,1.001407333333334)

# --- Row 15: Rotation-NoTilt (existing data, shifted down 4 rows) ---
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
Set-RowValues $ws 15 3 @(1,0.66,1.114299999999999,1,0.8,1.3,1.09,0.66,0.8871499999999994,0.9435749999999998,0.9940499999999998)

# --- Row 16: Rotation-60detTilt (existing data, shifted down 4 rows) ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
Set-RowValues $ws 16 3 @(0.9980463069183974,0.7996440047616012,1.063877739110392,0.9980463069183974,0.8815865060351997,1.170619099852799,1.048729296486398,0.7996440047616012,0.9317608719359964,0.9649035894271969,0.9937504921941311)

# --- Row 17: HexGrid-90degTilt5degRes (new row, existing data shifted down) ---
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
Set-RowValues $ws 17 3 @(0.9928149941888519,0.9949863574415481,0.9947052886848506,0.9928149941888519,0.9934784416931822,0.9955900074459716,0.9942882193454705,0.9949863574415481,0.9948458230631994,0.9938304086260256,0.9943105514666458)
Set-IndexCellStyle $ws 17

# --- Row 18: HexGrid-90degTilt22p5degRes (new row, existing data shifted down) ---
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
Set-RowValues $ws 18 3 @(0.999380597217319,1.002396640162602,0.9884982858496604,0.999380597217319,0.9990594747971785,0.9790972836567882,0.9919349066976585,1.002396640162602,0.9954474630061314,0.9974140301117251,0.9933945313968677)
Set-IndexCellStyle $ws 18

# --- Row 19: HexGrid-60degTilt5degRes (new row, existing data shifted down) ---
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
Set-RowValues $ws 19 3 @(0.9864979897053711,1.029237798317622,0.9850064625908554,0.9864979897053711,1.01581792504704,0.9645089886084285,0.9840139065222425,1.029237798317622,1.007122130454239,0.9968100600798049,0.9941805117985932)
Set-IndexCellStyle $ws 19
